# update to manual status column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# manualStatus column (I) now carries a bracketed text status ("[12]")
# instead of the raw numeric audit count, for rows 17 and 39.
$ws.Cells.Item(17, 9).Value = "[12]"
$ws.Cells.Item(39, 9).Value = "[12]"

# Widen the NOTES column (G) so manual notes are readable.
$ws.Columns.Item(7).ColumnWidth = 45.8333333333333

# Leave the selection where the user ended up after this edit.
$ws.Range("I40").Select()
